# Update the public EPEX Spot prices workbook:
#  - "Prix Spot" sheet: add a new day column BJ ("14-aug") with its 24 hourly prices
#  - "Gaz" sheet: append the new daily quote row (2025-08-12)
#  - "CO2" sheet: append the new daily quote row (2025-08-12)

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "Prix Spot": new column BJ (14-aug)
# ------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Prix Spot")

# Header cell - copy the formatting (bold font + border) from the previous
# header cell (BI1) so the new header matches the rest of the row.
$ws1.Range("BI1").Copy()
$ws1.Range("BJ1").PasteSpecial(-4122)
$ws1.Range("BJ1").Value = "14-aug"

# Hourly values for 14-aug
$ws1.Range("BJ2").Value = 98.09999999999999
$ws1.Range("BJ3").Value = 92
$ws1.Range("BJ4").Value = 88.45
$ws1.Range("BJ5").Value = 81.83
$ws1.Range("BJ6").Value = 82.70999999999999
$ws1.Range("BJ7").Value = 79.59
$ws1.Range("BJ8").Value = 84.89
$ws1.Range("BJ9").Value = 103.11
$ws1.Range("BJ10").Value = 99.31
$ws1.Range("BJ11").Value = 90.92
$ws1.Range("BJ12").Value = 76.66
$ws1.Range("BJ13").Value = 67.40000000000001
$ws1.Range("BJ14").Value = 54.01
$ws1.Range("BJ15").Value = 40.91
$ws1.Range("BJ16").Value = 54.9
$ws1.Range("BJ17").Value = 69.43000000000001
$ws1.Range("BJ18").Value = 78.98999999999999
$ws1.Range("BJ19").Value = 84.63
$ws1.Range("BJ20").Value = 102
$ws1.Range("BJ21").Value = 115.78
$ws1.Range("BJ22").Value = 137.99
$ws1.Range("BJ23").Value = 142.01
$ws1.Range("BJ24").Value = 120.75
$ws1.Range("BJ25").Value = 105.95

# ------------------------------------------------------------------
# Sheet "Gaz": append row 59 (2025-08-12)
# ------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Gaz")

# Force the date column to be stored as plain text (matching every other
# row in the column) instead of letting Excel auto-convert it to a date
# serial number.
$ws2.Range("A59").NumberFormat = "@"
$ws2.Range("A59").Value = "2025-08-12"
$ws2.Range("A59").Style = "Normal"
$ws2.Range("B59").Value = 31.225

# ------------------------------------------------------------------
# Sheet "CO2": append row 59 (2025-08-12)
# ------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("CO2")

$ws3.Range("A59").NumberFormat = "@"
$ws3.Range("A59").Value = "2025-08-12"
$ws3.Range("A59").Style = "Normal"
$ws3.Range("B59").Value = 70.84999999999999

$wb.Save()
